$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6, shifting the existing data rows (6-17) down to (7-18)
$ws.Rows("6:6").Insert()

# Populate the newly inserted row 6 with this week's new record
$ws.Range("A6").Value = 8
$ws.Range("B6").Value = "Terminal La Palmera de La Serena"
$ws.Range("C6").Value = "Coquimbo"
$ws.Range("D6").Value = 44414
$ws.Range("E6").Value = 4
$ws.Range("F6").Value = 100114007
$ws.Range("G6").Value = "Jengibre"
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 500
$ws.Range("K6").Value = 14000
$ws.Range("L6").Value = 15000
$ws.Range("M6").Value = 14500
$ws.Range("N6").Value = '$/caja 13 kilos'
$ws.Range("O6").Value = "Perú"
$ws.Range("P6").Value = 1115
$ws.Range("Q6").Value = 13
$ws.Range("R6").Value = "Hortaliza"
